$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - first sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 753
$ws1.Range("F9").Value = 416

# Sheet "全部类型" (all types) - fourth sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 753
$ws4.Range("F9").Value = 416
